$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.58102533333333
$ws.Range("H2").Value = 58.743076
$ws.Range("I2").Value = 0.3529199051285138
$ws.Range("J2").Value = 0.3529199051285138
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 37.407943255256
$ws.Range("R2").Value = 336.671489297304
$ws.Range("S2").Value = 0.00638475977750839
$ws.Range("T2").Value = 0.006384759777508391
$ws.Range("G3").Value = 19.58102533333333
$ws.Range("H3").Value = 58.743076
$ws.Range("I3").Value = 0.3529199051285138
$ws.Range("J3").Value = 0.3529199051285138
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 626.5000357686373
$ws.Range("R3").Value = 5638.500321917736
$ws.Range("S3").Value = 0.1069305575473235
$ws.Range("T3").Value = 0.1069305575473235
$ws.Range("G4").Value = 19.58102533333333
$ws.Range("H4").Value = 58.743076
$ws.Range("I4").Value = 0.3529199051285138
$ws.Range("J4").Value = 0.3529199051285138
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 741.312183368092
$ws.Range("R4").Value = 6671.809650312827
$ws.Range("S4").Value = 0.1265266090319065
$ws.Range("T4").Value = 0.1265266090319065
$ws.Range("G5").Value = 19.58102533333333
$ws.Range("H5").Value = 58.743076
$ws.Range("I5").Value = 0.3529199051285138
$ws.Range("J5").Value = 0.3529199051285138
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 662.5174259828383
$ws.Range("R5").Value = 5962.656833845544
$ws.Range("S5").Value = 0.1130779787717754
$ws.Range("T5").Value = 0.1130779787717754
$ws.Range("I6").Value = 0.2647616806631773
$ws.Range("J6").Value = 0.2647616806631773
$ws.Range("M6").Value = 1.910418
$ws.Range("N6").Value = 5.731254
$ws.Range("O6").Value = 0.01809124304049503
$ws.Range("P6").Value = 0.01809124304049503
$ws.Range("Q6").Value = 28.0635628155
$ws.Range("R6").Value = 252.5720653395
$ws.Range("S6").Value = 0.004789867912687473
$ws.Range("T6").Value = 0.004789867912687474
$ws.Range("I7").Value = 0.2647616806631773
$ws.Range("J7").Value = 0.2647616806631773
$ws.Range("O7").Value = 0.302988173785169
$ws.Range("P7").Value = 0.302988173785169
$ws.Range("S7").Value = 0.08021965811242818
$ws.Range("T7").Value = 0.0802196581124282
$ws.Range("I8").Value = 0.2647616806631773
$ws.Range("J8").Value = 0.2647616806631773
$ws.Range("M8").Value = 37.858701
$ws.Range("N8").Value = 113.576103
$ws.Range("O8").Value = 0.3585136661130873
$ws.Range("P8").Value = 0.3585136661130873
$ws.Range("Q8").Value = 556.1348530147499
$ws.Range("R8").Value = 5005.21367713275
$ws.Range("S8").Value = 0.0949206807808182
$ws.Range("T8").Value = 0.09492068078081821
$ws.Range("I9").Value = 0.2647616806631773
$ws.Range("J9").Value = 0.2647616806631773
$ws.Range("M9").Value = 33.83466466666667
$ws.Range("N9").Value = 101.503994
$ws.Range("O9").Value = 0.3204069170612486
$ws.Range("P9").Value = 0.3204069170612486
$ws.Range("Q9").Value = 497.0227652871666
$ws.Range("R9").Value = 4473.2048875845
$ws.Range("S9").Value = 0.08483147385724343
$ws.Range("T9").Value = 0.08483147385724345
$ws.Range("G10").Value = 19.14352733333333
$ws.Range("H10").Value = 57.430582
$ws.Range("I10").Value = 0.3450346309906436
$ws.Range("J10").Value = 0.3450346309906436
$ws.Range("M10").Value = 1.910418
$ws.Range("N10").Value = 5.731254
$ws.Range("O10").Value = 0.01809124304049503
$ws.Range("P10").Value = 0.01809124304049503
$ws.Range("Q10").Value = 36.572139201092
$ws.Range("R10").Value = 329.149252809828
$ws.Range("S10").Value = 0.006242105366639251
$ws.Range("T10").Value = 0.006242105366639251
$ws.Range("G11").Value = 19.14352733333333
$ws.Range("H11").Value = 57.430582
$ws.Range("I11").Value = 0.3450346309906436
$ws.Range("J11").Value = 0.3450346309906436
$ws.Range("O11").Value = 0.302988173785169
$ws.Range("P11").Value = 0.302988173785169
$ws.Range("Q11").Value = 612.5021726341614
$ws.Range("R11").Value = 5512.519553707452
$ws.Range("S11").Value = 0.1045414127364948
$ws.Range("T11").Value = 0.1045414127364948
$ws.Range("G12").Value = 19.14352733333333
$ws.Range("H12").Value = 57.430582
$ws.Range("I12").Value = 0.3450346309906436
$ws.Range("J12").Value = 0.3450346309906436
$ws.Range("M12").Value = 37.858701
$ws.Range("N12").Value = 113.576103
$ws.Range("O12").Value = 0.3585136661130873
$ws.Range("P12").Value = 0.3585136661130873
$ws.Range("Q12").Value = 724.749077397994
$ws.Range("R12").Value = 6522.741696581946
$ws.Range("S12").Value = 0.1236996304924319
$ws.Range("T12").Value = 0.1236996304924319
$ws.Range("G13").Value = 19.14352733333333
$ws.Range("H13").Value = 57.430582
$ws.Range("I13").Value = 0.3450346309906436
$ws.Range("J13").Value = 0.3450346309906436
$ws.Range("M13").Value = 33.83466466666667
$ws.Range("N13").Value = 101.503994
$ws.Range("O13").Value = 0.3204069170612486
$ws.Range("P13").Value = 0.3204069170612486
$ws.Range("Q13").Value = 647.714827860501
$ws.Range("R13").Value = 5829.433450744508
$ws.Range("S13").Value = 0.1105514823950777
$ws.Range("T13").Value = 0.1105514823950777
$ws.Range("G14").Value = 2.068613
$ws.Range("H14").Value = 6.205839
$ws.Range("I14").Value = 0.03728378321766519
$ws.Range("J14").Value = 0.0372837832176652
$ws.Range("M14").Value = 1.910418
$ws.Range("N14").Value = 5.731254
$ws.Range("O14").Value = 0.01809124304049503
$ws.Range("P14").Value = 0.01809124304049503
$ws.Range("Q14").Value = 3.951915510234
$ws.Range("R14").Value = 35.567239592106
$ws.Range("S14").Value = 0.0006745099836599107
$ws.Range("T14").Value = 0.0006745099836599108
$ws.Range("G15").Value = 2.068613
$ws.Range("H15").Value = 6.205839
$ws.Range("I15").Value = 0.03728378321766519
$ws.Range("J15").Value = 0.0372837832176652
$ws.Range("O15").Value = 0.302988173785169
$ws.Range("P15").Value = 0.302988173785169
$ws.Range("Q15").Value = 66.185814911606
$ws.Range("R15").Value = 595.672334204454
$ws.Range("S15").Value = 0.01129654538892251
$ws.Range("T15").Value = 0.01129654538892251
$ws.Range("G16").Value = 2.068613
$ws.Range("H16").Value = 6.205839
$ws.Range("I16").Value = 0.03728378321766519
$ws.Range("J16").Value = 0.0372837832176652
$ws.Range("M16").Value = 37.858701
$ws.Range("N16").Value = 113.576103
$ws.Range("O16").Value = 0.3585136661130873
$ws.Range("P16").Value = 0.3585136661130873
$ws.Range("Q16").Value = 78.315001051713
$ws.Range("R16").Value = 704.835009465417
$ws.Range("S16").Value = 0.01336674580793075
$ws.Range("T16").Value = 0.01336674580793075
$ws.Range("G17").Value = 2.068613
$ws.Range("H17").Value = 6.205839
$ws.Range("I17").Value = 0.03728378321766519
$ws.Range("J17").Value = 0.0372837832176652
$ws.Range("M17").Value = 33.83466466666667
$ws.Range("N17").Value = 101.503994
$ws.Range("O17").Value = 0.3204069170612486
$ws.Range("P17").Value = 0.3204069170612486
$ws.Range("Q17").Value = 69.99082718010735
$ws.Range("R17").Value = 629.9174446209661
$ws.Range("S17").Value = 0.01194598203715202
$ws.Range("T17").Value = 0.01194598203715203

Write-Host "Updated 174 cells"
